$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking strings so Excel does not
# silently convert them to Number cells (matches the source data,
# which stores these as plain text / inline strings).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "60.848.32"
$ws.Range("E2").Value = "  +0.48%  "
# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.592.74"
$ws.Range("E3").Value = "  +0.31%  "
# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "
# Row 5 - BNB
Set-TextValue $ws.Range("D5") "522.65"
$ws.Range("E5").Value = "  +3.07%  "
# Row 6 - Solana
Set-TextValue $ws.Range("D6") "154.20"
$ws.Range("E6").Value = "  +0.68%  "
# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.06%  "
# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.594"
$ws.Range("E8").Value = "  +2.76%  "
# Row 9 - Toncoin
Set-TextValue $ws.Range("D9") "6.71"
$ws.Range("E9").Value = "  +2.02%  "
# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.51%  "
# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.347"
$ws.Range("E11").Value = "  +0.01%  "
# Row 12 - TRON
$ws.Range("E12").Value = "  +1.39%  "
# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "3.048.70"
$ws.Range("E13").Value = "  +0.23%  "
# Row 14 - WrappedBTC
Set-TextValue $ws.Range("D14") "60.863.91"
$ws.Range("E14").Value = "  +0.59%  "
# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "21.62"
$ws.Range("E15").Value = "  +0.40%  "
# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.18%  "
# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.597.12"
$ws.Range("E17").Value = "  +0.22%  "
# Row 18 - Polkadot
$ws.Range("E18").Value = "  -0.87%  "
# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "352.79"
$ws.Range("E19").Value = "  +2.14%  "
# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "10.57"
$ws.Range("E20").Value = "  +1.37%  "
# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.22"
$ws.Range("E21").Value = "  +1.66%  "
# Row 22 - Dai
$ws.Range("E22").Value = "  +0.27%  "
# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "60.80"
$ws.Range("E23").Value = "  +1.36%  "
# Row 24 - Polygon
Set-TextValue $ws.Range("D24") "0.426"
$ws.Range("E24").Value = "  +1.52%  "
# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.71%  "
# Row 26 - WrappedeETH
Set-TextValue $ws.Range("D26") "2.711.35"
$ws.Range("E26").Value = "  +0.35%  "
# Row 27 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D27") "0.999"
$ws.Range("E27").Value = "  +0.14%  "
# Row 28 - PEPE
Set-TextValue $ws.Range("D28") "0.0₃0843"
$ws.Range("E28").Value = "  -0.19%  "
# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "7.36"
$ws.Range("E29").Value = "  +0.16%  "
# Row 30 - USDe
$ws.Range("E30").Value = "  -0.05%  "
# Row 31 - Aptos
Set-TextValue $ws.Range("D31") "6.33"
$ws.Range("E31").Value = "  +10.76%  "
# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "19.36"
$ws.Range("E32").Value = "  +0.14%  "
# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +2.85%  "
# Row 34 - Monero
Set-TextValue $ws.Range("D34") "148.19"
$ws.Range("E34").Value = "  -3.46%  "
# Row 35 - NEARProtocol
Set-TextValue $ws.Range("D35") "4.15"
$ws.Range("E35").Value = "  +4.20%  "
# Row 36 - SuiNetwork
Set-TextValue $ws.Range("D36") "0.934"
$ws.Range("E36").Value = "  +8.97%  "
# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +0.75%  "
# Row 38 - Stacks
$ws.Range("E38").Value = "  +1.80%  "
# Row 39 - Filecoin/Fetch.AI swap
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D39") "0.849"
$ws.Range("E39").Value = "  -0.31%  "
# Row 40 - Fetch.AI/Filecoin swap
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D40") "3.79"
$ws.Range("E40").Value = "  +1.03%  "
# Row 41 - OKB
$ws.Range("E41").Value = "  +1.45%  "
# Row 42 - Bittensor
Set-TextValue $ws.Range("D42") "287.58"
$ws.Range("E42").Value = "  -2.74%  "
# Row 43 - Stellar
$ws.Range("E43").Value = "  +1.76%  "
# Row 44 - Mantle
$ws.Range("E44").Value = "  +0.58%  "
# Row 45 - Hedera
$ws.Range("E45").Value = "  +0.59%  "
# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.07%  "
# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "19.56"
$ws.Range("E47").Value = "  -1.44%  "
# Row 48 - VeChain
$ws.Range("E48").Value = "  +2.22%  "
# Row 49 - RenderToken
$ws.Range("E49").Value = "  +0.47%  "
# Row 50 - WhiteBITCoin
Set-TextValue $ws.Range("D50") "10.32"
$ws.Range("E50").Value = "  +0.14%  "
# Row 51 - InjectiveProtocol
Set-TextValue $ws.Range("D51") "19.07"
$ws.Range("E51").Value = "  +8.46%  "
